$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell D4: "a" (least_popular row, custom scoring column)
$ws.Range("D4").Value = "a"

# A5 becomes text "1.2" instead of numeric 1.2 (cell keeps its existing text number format)
$ws.Range("A5").Value = "1.2"

# Row 9: aggregate_difficulty answer changes from B to A
$ws.Range("C9").Value = "A"

# A12 becomes text "2.4" instead of numeric 2.4
$ws.Range("A12").Value = "2.4"

# Row 18 becomes question id "2.10" (text) with answer B (shift of the whole block below)
$ws.Range("A18").Value = "2.10"
$ws.Range("C18").Value = "B"

# Rows 19-23 shift up by one question id, and answers are updated.
# These A cells carry a "@" (text) number format, so writing a numeric
# Value directly would coerce it to a text cell. Temporarily switch the
# format to General while assigning the number, then restore the
# original "@" format via a format-only paste (keeps the original style
# index instead of minting a new one).
$ws.Range("A19").NumberFormat = "General"
$ws.Range("A19").Value = 3.1
$ws.Range("C19").Value = "C"

$ws.Range("A20").NumberFormat = "General"
$ws.Range("A20").Value = 3.2
$ws.Range("C20").Value = "B"

$ws.Range("A21").NumberFormat = "General"
$ws.Range("A21").Value = 3.3
$ws.Range("C21").Value = "A"

$ws.Range("A22").NumberFormat = "General"
$ws.Range("A22").Value = 3.4
$ws.Range("C22").Value = "D"

$ws.Range("A23").NumberFormat = "General"
$ws.Range("A23").Value = 3.5
$ws.Range("C23").Value = "D"

# Row 24: new row, question id 3.6, round_type aggregate_difficulty, answer A
$ws.Range("A24").Value = 3.6
$ws.Range("B24").Value = "aggregate_difficulty"
$ws.Range("C24").Value = "A"

# Restore the original "@" text format for column A (rows 19-24) by
# copying the format (only) from row 18's A cell, which still carries
# the original style.
$ws.Range("A18").Copy()
$ws.Range("A19:A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selected cell in the sheet view
$ws.Range("F5").Select()

# Nudge the sheet's default column width to match the updated layout
$ws.StandardWidth = 8.47265625
